# updated prod smoke test browsers
# Applies the 2020 -> 2021 MSRP refresh for RC 300 / RC 350 rows, plus the
# new "Black Line" special-edition trims appended at the bottom of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Year / Base MSRP updates for existing RC 300 / RC 350 rows (2020 -> 2021) ---
# row, year, new base MSRP
$updates = @(
    @(2,  2021, 42120),
    @(3,  2021, 46590),
    @(4,  2021, 44810),
    @(5,  2021, 48765),
    @(6,  2021, 45050),
    @(7,  2021, 49520),
    @(8,  2021, 47215),
    @(9,  2021, 51130),
    @(53, 2021, 65875),
    @(54, 2021, 96675)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 3).Value = $u[1]
    $ws.Cells.Item($row, 4).Value = $u[2]
}

# --- Append the four new Black Line special-edition trim rows (95-98) ---
# Shared strings are written in two passes (all trim codes, then all trim
# names) so the new entries land in the sharedStrings table in the same
# order Excel produced them in the source workbook.
$trimCodes = @("9203SE", "9207SE", "9213SE", "9217SE")
$trimNames = @("RC 300 F SPORT Black Line", "RC 300 AWD F SPORT Black Line", "RC 350 F SPORT Black Line", "RC 350 AWD F SPORT Black Line")
$years     = @(2021, 2021, 2021, 2021)
$msrps     = @(48735, 50910, 51665, 53275)

for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item(95 + $i, 1).Value = $trimCodes[$i]
}
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item(95 + $i, 2).Value = $trimNames[$i]
}

for ($i = 0; $i -lt 4; $i++) {
    $rowNum = 95 + $i
    $ws.Cells.Item($rowNum, 3).Value = $years[$i]
    $ws.Cells.Item($rowNum, 4).Value = $msrps[$i]
    $ws.Cells.Item($rowNum, 4).NumberFormat = '"$"#,##0_);[Red]("$"#,##0)'
    $ws.Cells.Item($rowNum, 5).Value = 1025
    $ws.Cells.Item($rowNum, 5).NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'
}

# --- Restore the view state captured in the saved workbook ---
$ws.Range("C55").Select()
